{"js": "// Update page numbers / section references in the List of Figures, List of\n// Tables, and Table of Contents tables to match the corrected pagination of\n// the source PDF.\n//\n// Each entry is [tableIndex, rowIndex, columnIndex, expectedOldText, newText]\n// (all indices are 0-based, matching Office.js Table.getCell semantics).\n// tableIndex 0 = \"List of Figures\", 1 = \"List of Tables\", 2 = \"Table of Contents\".\nconst edits = [\n  [0, 12, 2, \"20\", \"21\"],\n  [0, 14, 2, \"21\", \"22\"],\n  [0, 15, 1, \"System Workflows\", \"System Workflows (Section 2.4)\"],\n  [0, 16, 2, \"21\", \"22\"],\n  [0, 17, 2, \"21\", \"23\"],\n  [0, 18, 1, \"Technical Stack\", \"Technical Stack (Section 3.1)\"],\n  [0, 19, 2, \"22\", \"23\"],\n  [0, 21, 2, \"23\", \"25\"],\n  [0, 22, 2, \"23\", \"25\"],\n  [0, 23, 2, \"24\", \"26\"],\n  [0, 25, 2, \"24\", \"26\"],\n  [0, 26, 2, \"26\", \"28\"],\n  [0, 27, 2, \"27\", \"29\"],\n  [0, 28, 2, \"28\", \"30\"],\n  [0, 30, 2, \"29\", \"32\"],\n  [0, 31, 2, \"29\", \"33\"],\n  [0, 32, 2, \"30\", \"34\"],\n  [0, 33, 2, \"30\", \"34\"],\n  [1, 3, 2, \"22\", \"24\"],\n  [1, 4, 2, \"25\", \"27\"],\n  [1, 5, 2, \"31\", \"35\"],\n  [1, 6, 2, \"32\", \"36\"],\n  [1, 7, 2, \"33\", \"37\"],\n  [1, 8, 2, \"33\", \"37\"],\n  [2, 11, 2, \"21\", \"22\"],\n  [2, 12, 2, \"22\", \"23\"],\n  [2, 13, 2, \"22\", \"23\"],\n  [2, 14, 2, \"23\", \"25\"],\n  [2, 15, 2, \"24\", \"26\"],\n  [2, 16, 2, \"25\", \"28\"],\n  [2, 17, 2, \"27\", \"29\"],\n  [2, 18, 2, \"28\", \"30\"],\n  [2, 19, 2, \"29\", \"32\"],\n  [2, 20, 2, \"29\", \"32\"],\n  [2, 21, 2, \"31\", \"35\"],\n  [2, 22, 2, \"31\", \"35\"],\n  [2, 23, 2, \"31\", \"35\"],\n  [2, 24, 2, \"32\", \"36\"],\n  [2, 25, 2, \"33\", \"37\"],\n  [2, 26, 2, \"33\", \"37\"],\n  [2, 27, 2, \"33\", \"38\"],\n  [2, 28, 2, \"34\", \"39\"],\n  [2, 29, 2, \"34\", \"39\"],\n  [2, 30, 2, \"34\", \"39\"],\n  [2, 31, 2, \"34\", \"39\"],\n  [2, 32, 2, \"35\", \"40\"],\n  [2, 33, 2, \"36\", \"41\"],\n  [2, 34, 2, \"36\", \"41\"],\n  [2, 35, 2, \"36\", \"41\"],\n  [2, 36, 2, \"36\", \"41\"],\n  [2, 37, 2, \"37\", \"42\"]\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nfor (const [tableIndex, rowIndex, colIndex, oldText, newText] of edits) {\n  const table = tables.items[tableIndex];\n  const cell = table.getCell(rowIndex, colIndex);\n  const body = cell.body;\n  body.load(\"paragraphs\");\n  await context.sync();\n\n  const paragraph = body.paragraphs.items[0];\n  paragraph.load(\"text\");\n  await context.sync();\n\n  // Defensive check: only replace when the current text matches what the\n  // diff expects to be there, so we never clobber the wrong cell.\n  const currentText = paragraph.text.trim();\n  if (currentText !== oldText) {\n    throw new Error(\n      `Unexpected text in table ${tableIndex} row ${rowIndex} col ${colIndex}: ` +\n      `expected \"${oldText}\" but found \"${currentText}\"`\n    );\n  }\n\n  const range = paragraph.getRange();\n  range.insertText(newText, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Update page numbers / section references in the List of Figures, List of\n# Tables, and Table of Contents tables to match the corrected pagination of\n# the source PDF.\n#\n# Each entry is (tableIndex, rowIndex, columnIndex, expectedOldText, newText)\n# using 1-based indices, matching the Word COM object model (Document.Tables,\n# Table.Cell(row, col)). Table 1 = \"List of Figures\", Table 2 = \"List of\n# Tables\", Table 3 = \"Table of Contents\".\n$edits = @(\n    @(1, 13, 3, \"20\", \"21\"),\n    @(1, 15, 3, \"21\", \"22\"),\n    @(1, 16, 2, \"System Workflows\", \"System Workflows (Section 2.4)\"),\n    @(1, 17, 3, \"21\", \"22\"),\n    @(1, 18, 3, \"21\", \"23\"),\n    @(1, 19, 2, \"Technical Stack\", \"Technical Stack (Section 3.1)\"),\n    @(1, 20, 3, \"22\", \"23\"),\n    @(1, 22, 3, \"23\", \"25\"),\n    @(1, 23, 3, \"23\", \"25\"),\n    @(1, 24, 3, \"24\", \"26\"),\n    @(1, 26, 3, \"24\", \"26\"),\n    @(1, 27, 3, \"26\", \"28\"),\n    @(1, 28, 3, \"27\", \"29\"),\n    @(1, 29, 3, \"28\", \"30\"),\n    @(1, 31, 3, \"29\", \"32\"),\n    @(1, 32, 3, \"29\", \"33\"),\n    @(1, 33, 3, \"30\", \"34\"),\n    @(1, 34, 3, \"30\", \"34\"),\n    @(2, 4, 3, \"22\", \"24\"),\n    @(2, 5, 3, \"25\", \"27\"),\n    @(2, 6, 3, \"31\", \"35\"),\n    @(2, 7, 3, \"32\", \"36\"),\n    @(2, 8, 3, \"33\", \"37\"),\n    @(2, 9, 3, \"33\", \"37\"),\n    @(3, 12, 3, \"21\", \"22\"),\n    @(3, 13, 3, \"22\", \"23\"),\n    @(3, 14, 3, \"22\", \"23\"),\n    @(3, 15, 3, \"23\", \"25\"),\n    @(3, 16, 3, \"24\", \"26\"),\n    @(3, 17, 3, \"25\", \"28\"),\n    @(3, 18, 3, \"27\", \"29\"),\n    @(3, 19, 3, \"28\", \"30\"),\n    @(3, 20, 3, \"29\", \"32\"),\n    @(3, 21, 3, \"29\", \"32\"),\n    @(3, 22, 3, \"31\", \"35\"),\n    @(3, 23, 3, \"31\", \"35\"),\n    @(3, 24, 3, \"31\", \"35\"),\n    @(3, 25, 3, \"32\", \"36\"),\n    @(3, 26, 3, \"33\", \"37\"),\n    @(3, 27, 3, \"33\", \"37\"),\n    @(3, 28, 3, \"33\", \"38\"),\n    @(3, 29, 3, \"34\", \"39\"),\n    @(3, 30, 3, \"34\", \"39\"),\n    @(3, 31, 3, \"34\", \"39\"),\n    @(3, 32, 3, \"34\", \"39\"),\n    @(3, 33, 3, \"35\", \"40\"),\n    @(3, 34, 3, \"36\", \"41\"),\n    @(3, 35, 3, \"36\", \"41\"),\n    @(3, 36, 3, \"36\", \"41\"),\n    @(3, 37, 3, \"36\", \"41\"),\n    @(3, 38, 3, \"37\", \"42\")\n)\n\n$d = $word.ActiveDocument\n\nforeach ($edit in $edits) {\n    $tableIndex  = $edit[0]\n    $rowIndex    = $edit[1]\n    $colIndex    = $edit[2]\n    $oldText     = $edit[3]\n    $newText     = $edit[4]\n\n    $table = $d.Tables.Item($tableIndex)\n    $cell = $table.Cell($rowIndex, $colIndex)\n\n    # Cell.Range.Text includes the trailing cell-mark characters (CR + BEL);\n    # strip them before comparing against the expected old value.\n    $currentText = $cell.Range.Text.TrimEnd([char]13, [char]7)\n\n    if ($currentText -ne $oldText) {\n        throw \"Unexpected text in table $tableIndex row $rowIndex col $colIndex - expected '$oldText' but found '$currentText'\"\n    }\n\n    $cell.Range.Text = $newText\n}\n"}
